$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# player_x / camera_x (column A / E): 101 -> 79 (and its negative)
$ws.Range("A2").Value = 79
$ws.Range("E2").Value = -79

# player_y / camera_y (column B / F): 78 -> 74 (and its negative)
$ws.Range("B2").Value = 74
$ws.Range("F2").Value = -74
